# Fix permissions to author in home
# - Correct missing accents in the header row (shared strings):
#     "Correo electronico (*)" -> "Correo electrónico (*)"
#     "Subarea 1 (*)"          -> "Subárea 1 (*)"
#     "Subarea 2"              -> "Subárea 2"
#     "Subarea 3"              -> "Subárea 3"
# - Move the sheet view / selection from A3 to Q1 (scrolled so column J is
#   the left-most visible column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("D1").Value = "Correo electrónico (*)"
$ws.Range("F1").Value = "Subárea 1 (*)"
$ws.Range("G1").Value = "Subárea 2"
$ws.Range("H1").Value = "Subárea 3"

# Scroll the window so column J is left-most, then move the selection to Q1
# (mirrors topLeftCell="J1" / selection activeCell="Q1" sqref="Q1").
$win = $excel.ActiveWindow
$win.DisplayGridlines = $true
$win.ScrollColumn = 10
$win.ScrollRow = 1
$ws.Range("Q1").Select()
